$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 2599
$ws.Range("E2").Value = 77
$ws.Range("F2").Value = 77
$ws.Range("G2").Value = 26
$ws.Range("H2").Value = 11
$ws.Range("I2").Value = 14
$ws.Range("J2").Value = -3
$ws.Range("K2").Value = 2571
$ws.Range("L2").Value = 1608
$ws.Range("M2").Value = 963
$ws.Range("N2").Value = 1040
$ws.Range("O2").Value = -76
$ws.Range("P2").Value = 64
$ws.Range("Q2").Value = -244
$ws.Range("R2").Value = -24
$ws.Range("S2").Value = 245
$ws.Range("T2").Value = 21
$ws.Range("U2").Value = -265
$ws.Range("V2").Value = 1130
$ws.Range("W2").Value = 2.97
$ws.Range("X2").Value = 0.42
$ws.Range("Y2").Value = 1.39
$ws.Range("Z2").Value = 0.45
$ws.Range("AA2").Value = 166.88
$ws.Range("AB2").Value = 1642.13
$ws.Range("AC2").Value = 113
$ws.Range("AD2").Value = 29.26
$ws.Range("AE2").Value = 9233
$ws.Range("AF2").Value = 0.36
$ws.Range("AG2").Value = 50
$ws.Range("AI2").Value = 39.03
$ws.Range("AJ2").Value = 12712747

$ws.Range("D3").Value = 2043
$ws.Range("E3").Value = 47
$ws.Range("F3").Value = 26
$ws.Range("G3").Value = -8
$ws.Range("H3").Value = -45
$ws.Range("I3").Value = -33
$ws.Range("J3").Value = -12
$ws.Range("K3").Value = 2366
$ws.Range("L3").Value = 1454
$ws.Range("M3").Value = 913
$ws.Range("N3").Value = 1001
$ws.Range("O3").Value = -88
$ws.Range("P3").Value = 64
$ws.Range("Q3").Value = 306
$ws.Range("R3").Value = -5
$ws.Range("S3").Value = -236
$ws.Range("T3").Value = 29
$ws.Range("U3").Value = 276
$ws.Range("V3").Value = 916
$ws.Range("W3").Value = 2.3
$ws.Range("X3").Value = -2.2
$ws.Range("Y3").Value = -3.23
$ws.Range("Z3").Value = -1.82
$ws.Range("AA3").Value = 159.32
$ws.Range("AB3").Value = 1580.75
$ws.Range("AC3").Value = -260
$ws.Range("AD3").Value = -12.1
$ws.Range("AE3").Value = 8888
$ws.Range("AF3").Value = 0.35
$ws.Range("AG3").Value = 30
$ws.Range("AH3").Value = 0.96
$ws.Range("AI3").Value = -10.24
$ws.Range("AJ3").Value = 12712747

$ws.Range("D4").Value = 2110
$ws.Range("E4").Value = 52
$ws.Range("F4").Value = 52
$ws.Range("G4").Value = 19
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = -6
$ws.Range("K4").Value = 2346
$ws.Range("L4").Value = 1432
$ws.Range("M4").Value = 913
$ws.Range("N4").Value = 1007
$ws.Range("O4").Value = -94
$ws.Range("P4").Value = 64
$ws.Range("Q4").Value = 22
$ws.Range("R4").Value = -95
$ws.Range("S4").Value = 15
$ws.Range("T4").Value = 48
$ws.Range("U4").Value = -26
$ws.Range("V4").Value = 942
$ws.Range("W4").Value = 2.45
$ws.Range("X4").Value = 0.14
$ws.Range("Y4").Value = 0.94
$ws.Range("Z4").Value = 0.13
$ws.Range("AA4").Value = 156.88
$ws.Range("AB4").Value = 1590.1
$ws.Range("AC4").Value = 74
$ws.Range("AD4").Value = 55.89
$ws.Range("AE4").Value = 8944
$ws.Range("AF4").Value = 0.46
$ws.Range("AG4").Value = 40
$ws.Range("AH4").Value = 0.97
$ws.Range("AI4").Value = 47.77
$ws.Range("AJ4").Value = 12712747

$ws.Range("D5").Value = 2375
$ws.Range("E5").Value = 25
$ws.Range("F5").Value = 25
$ws.Range("G5").Value = 15
$ws.Range("H5").Value = 13
$ws.Range("I5").Value = 20
$ws.Range("J5").Value = -7
$ws.Range("K5").Value = 2280
$ws.Range("L5").Value = 1356
$ws.Range("M5").Value = 923
$ws.Range("N5").Value = 1005
$ws.Range("O5").Value = -82
$ws.Range("P5").Value = 64
$ws.Range("Q5").Value = -210
$ws.Range("R5").Value = 254
$ws.Range("S5").Value = 6
$ws.Range("T5").Value = 7
$ws.Range("U5").Value = -217
$ws.Range("V5").Value = 933
$ws.Range("W5").Value = 1.06
$ws.Range("X5").Value = 0.54
$ws.Range("Y5").Value = 1.99
$ws.Range("Z5").Value = 0.55
$ws.Range("AA5").Value = 146.96
$ws.Range("AB5").Value = 1586.63
$ws.Range("AC5").Value = 157
$ws.Range("AD5").Value = 22.1
$ws.Range("AE5").Value = 8925
$ws.Range("AF5").Value = 0.39
$ws.Range("AG5").Value = 30
$ws.Range("AH5").Value = 0.86
$ws.Range("AI5").Value = 16.9
$ws.Range("AJ5").Value = 12712747

$ws.Range("D6").Value = 2167
$ws.Range("E6").Value = 62
$ws.Range("F6").Value = 62
$ws.Range("G6").Value = 6
$ws.Range("H6").Value = -7
$ws.Range("I6").Value = -3
$ws.Range("K6").Value = 2284
$ws.Range("L6").Value = 1367
$ws.Range("M6").Value = 917
$ws.Range("N6").Value = 999
$ws.Range("P6").Value = 64
$ws.Range("Q6").Value = 30
$ws.Range("R6").Value = -16
$ws.Range("S6").Value = -45
$ws.Range("T6").Value = 8
$ws.Range("U6").Value = 22
$ws.Range("V6").Value = 888
$ws.Range("W6").Value = 2.85
$ws.Range("X6").Value = -0.35
$ws.Range("Y6").Value = -0.3
$ws.Range("Z6").Value = -0.33
$ws.Range("AA6").Value = 149.18
$ws.Range("AB6").Value = 1579.26
$ws.Range("AC6").Value = -23
$ws.Range("AD6").Value = -226.79
$ws.Range("AE6").Value = 8870
$ws.Range("AF6").Value = 0.6
$ws.Range("AG6").Value = 30
$ws.Range("AH6").Value = 0.56
$ws.Range("AI6").Value = -113.48
$ws.Range("AJ6").Value = 12712747

$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()
